$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The underlying TPM data was recomputed. The sending/target cluster
# combination that used "MuSCs" as a target cluster is no longer present,
# so the last two rows of the old table (rows 8 and 9) are dropped and the
# remaining rows are refreshed with the newly calculated values.
$ws.Range("A8:T9").Delete() | Out-Null

# Row 2: MuSCs -> Il27/Il27ra -> ECs
$ws.Range("A2").Value = "MuSCs"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.1055543333333333
$ws.Range("H2").Value = 0.316663
$ws.Range("I2").Value = 0.1206126635621877
$ws.Range("J2").Value = 0.1206126635621877
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 1.740493666666667
$ws.Range("N2").Value = 5.221481
$ws.Range("O2").Value = 0.5281470583624094
$ws.Range("P2").Value = 0.5281470583624094
$ws.Range("Q2").Value = 0.1837166486558889
$ws.Range("R2").Value = 1.653449837903
$ws.Range("S2").Value = 0.06370122346162439
$ws.Range("T2").Value = 0.06370122346162439

# Row 3: MuSCs -> Il27/Il27ra -> FAPs
$ws.Range("A3").Value = "MuSCs"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.1055543333333333
$ws.Range("H3").Value = 0.316663
$ws.Range("I3").Value = 0.1206126635621877
$ws.Range("J3").Value = 0.1206126635621877
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 1.169098333333333
$ws.Range("N3").Value = 3.507295
$ws.Range("O3").Value = 0.354759030447336
$ws.Range("P3").Value = 0.3547590304473361
$ws.Range("Q3").Value = 0.1234033951761111
$ws.Range("R3").Value = 1.110630556585
$ws.Range("S3").Value = 0.04278843158499244
$ws.Range("T3").Value = 0.04278843158499245

# Row 4: MuSCs -> Il27/Il27ra -> Resolving-Mac (was "-> MuSCs" before)
$ws.Range("A4").Value = "MuSCs"
$ws.Range("D4").Value = "Resolving-Mac"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.1055543333333333
$ws.Range("H4").Value = 0.316663
$ws.Range("I4").Value = 0.1206126635621877
$ws.Range("J4").Value = 0.1206126635621877
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.3858796666666667
$ws.Range("N4").Value = 1.157639
$ws.Range("O4").Value = 0.1170939111902545
$ws.Range("P4").Value = 0.1170939111902545
$ws.Range("Q4").Value = 0.04073127096188889
$ws.Range("R4").Value = 0.366581438657
$ws.Range("S4").Value = 0.01412300851557085
$ws.Range("T4").Value = 0.01412300851557085

# Row 5: Resolving-Mac -> Il27/Il27ra -> ECs (was "MuSCs -> Resolving-Mac" before)
$ws.Range("A5").Value = "Resolving-Mac"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.7695970000000001
$ws.Range("H5").Value = 2.308791
$ws.Range("I5").Value = 0.8793873364378123
$ws.Range("J5").Value = 0.8793873364378122
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 1.740493666666667
$ws.Range("N5").Value = 5.221481
$ws.Range("O5").Value = 0.5281470583624094
$ws.Range("P5").Value = 0.5281470583624094
$ws.Range("Q5").Value = 1.339478704385667
$ws.Range("R5").Value = 12.055308339471
$ws.Range("S5").Value = 0.464445834900785
$ws.Range("T5").Value = 0.4644458349007849

# Row 6: Resolving-Mac -> Il27/Il27ra -> FAPs
$ws.Range("A6").Value = "Resolving-Mac"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 0.7695970000000001
$ws.Range("H6").Value = 2.308791
$ws.Range("I6").Value = 0.8793873364378123
$ws.Range("J6").Value = 0.8793873364378122
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 1.169098333333333
$ws.Range("N6").Value = 3.507295
$ws.Range("O6").Value = 0.354759030447336
$ws.Range("P6").Value = 0.3547590304473361
$ws.Range("Q6").Value = 0.8997345700383335
$ws.Range("R6").Value = 8.097611130345001
$ws.Range("S6").Value = 0.3119705988623436
$ws.Range("T6").Value = 0.3119705988623436

# Row 7: Resolving-Mac -> Il27/Il27ra -> Resolving-Mac (was "-> FAPs" before)
$ws.Range("A7").Value = "Resolving-Mac"
$ws.Range("D7").Value = "Resolving-Mac"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 0.7695970000000001
$ws.Range("H7").Value = 2.308791
$ws.Range("I7").Value = 0.8793873364378123
$ws.Range("J7").Value = 0.8793873364378122
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.3858796666666667
$ws.Range("N7").Value = 1.157639
$ws.Range("O7").Value = 0.1170939111902545
$ws.Range("P7").Value = 0.1170939111902545
$ws.Range("Q7").Value = 0.2969718338276667
$ws.Range("R7").Value = 2.672746504449
$ws.Range("S7").Value = 0.1029709026746836
$ws.Range("T7").Value = 0.1029709026746836
